# Updates "Estado de Cuenta" detail rows (16-22): the previous mora-period
# entries are replaced with the new ones (Periodo Mora / Valor Mora /
# Salario Basico), per "Elimina EC anteriores y se agregan nuevos,
# se modifica base de datos".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: Periodo 1811 -> 1901, Valor Mora 15625 -> 31249, Salario Basico 3000000 -> 1200000
$ws.Range("E16").Value = "1901"
$ws.Range("F16").Value = 31249
$ws.Range("G16").Value = 1200000

# Row 17: Periodo 1812 stays, Valor Mora 31249 stays, Salario Basico 3000000 -> 1200000
$ws.Range("E17").Value = "1812"
$ws.Range("F17").Value = 31249
$ws.Range("G17").Value = 1200000

# Row 18: Periodo 1901 -> 1811, Valor Mora 31249 -> 15625, Salario Basico 3000000 -> 1200000
$ws.Range("E18").Value = "1811"
$ws.Range("F18").Value = 15625
$ws.Range("G18").Value = 1200000

# Row 19: Periodo 1902 -> 1905, Valor Mora 20800 -> 52000, Salario Basico 1300000 stays
$ws.Range("E19").Value = "1905"
$ws.Range("F19").Value = 52000
$ws.Range("G19").Value = 1300000

# Row 20: Periodo 1903 -> 1904, Valor Mora 52000 stays, Salario Basico 1300000 stays
$ws.Range("E20").Value = "1904"
$ws.Range("F20").Value = 52000
$ws.Range("G20").Value = 1300000

# Row 21: Periodo 1904 -> 1903, Valor Mora 52000 stays, Salario Basico 1300000 stays
$ws.Range("E21").Value = "1903"
$ws.Range("F21").Value = 52000
$ws.Range("G21").Value = 1300000

# Row 22: Periodo 1905 -> 1902, Valor Mora 52000 -> 20800, Salario Basico 1300000 stays
$ws.Range("E22").Value = "1902"
$ws.Range("F22").Value = 20800
$ws.Range("G22").Value = 1300000
